# Resize the columns on the "Data" sheet.
# The stored OOXML <col width="..."> value is the Excel ColumnWidth
# plus a constant offset of 5/6 (≈0.8333333) character units, so we
# back that offset out of each target stored width before assigning
# ColumnWidth via COM.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$offset = 5 / 6

# Target stored widths (what should end up in the XML) per column A..J
$targetWidths = @(15, 9, 7, 8, 8, 7, 5, 13, 6, 50)

for ($i = 0; $i -lt $targetWidths.Length; $i++) {
    $col = $i + 1
    $ws.Columns.Item($col).ColumnWidth = $targetWidths[$i] - $offset
}
